# Add "Abyss Actor - Extra" and related new card-id rows to the SPDS-JP sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("SPDS-JP")

# --- Fill in card names that were missing on existing rows, and append new
#     rows 33-37, in the exact order the new shared strings were authored so
#     that the shared-string table indices line up with the source edit. ---

# 1) "Abyss Actor - Extra" card id 100405020
$ws.Range("A21").Value = "Abyss Actor - Extra"

# 2) Append new rows 33-37 following the existing pattern (B=id, C=":", E=",")
$ws.Range("B33").Value = 100405032
$ws.Range("C33").Value = ":"
$ws.Range("E33").Value = ","

$ws.Range("B34").Value = 100405033
$ws.Range("C34").Value = ":"
$ws.Range("E34").Value = ","

$ws.Range("A35").Value = "Casting out the Darklords"
$ws.Range("B35").Value = 100405034
$ws.Range("C35").Value = ":"
$ws.Range("E35").Value = ","

$ws.Range("B36").Value = 100405035
$ws.Range("C36").Value = ":"
$ws.Range("E36").Value = ","

$ws.Range("A37").Value = "Darklords Falling from Grace"
$ws.Range("B37").Value = 100405036
$ws.Range("C37").Value = ":"
$ws.Range("E37").Value = ","

# 3) Fill in remaining card names that were missing on existing rows
$ws.Range("A23").Value = "Abyss Actor - Wild Hope"
$ws.Range("A24").Value = "Abyss Script - Fantasy Magic"
$ws.Range("A28").Value = "Abyss Script - Rise of the Dark Ruler"

# --- Update the selection to match where the user ended up after editing ---
$ws.Range("D36").Select()
